# Edit script: version 0.06, store team data
# Applies two structural edits to the project log docx:
#  1) Remove the paragraph-mark run properties (rFonts hint=eastAsia) from
#     the "完成功能：调试前后比赛日获取函数" paragraph.
#  2) Remove the paragraph-mark run properties from the last paragraph
#     ("下一步工作：...") and restructure/extend its content: the trailing
#     "尝试将数据存入sql数据库..." runs move to immediately follow "工作："
#     (ahead of the _GoBack bookmark), and five new paragraphs are appended
#     describing the 17/6/24 entry (done / resolved / unresolved / next
#     steps), ending with a fresh "下一步工作：...存入数据库。" paragraph that
#     now carries the _GoBack bookmark.

$d = $word.ActiveDocument

# Locate the two target paragraphs by their distinctive text so the script
# is resilient to any incidental paragraph-numbering differences.
$targetCount = $d.Paragraphs.Count

$idx40 = 0
$idx42 = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text
    if ($ptext.StartsWith("完成功能") -and $ptext.Contains("调试前后比赛日")) {
        $idx40 = $i
    }
    if ($ptext.StartsWith("下一步") -and $ptext.Contains("尝试") -and $ptext.Contains("数据库以供网页和远程应用访问")) {
        $idx42 = $i
    }
}

if ($idx40 -eq 0) { throw "could not find the '完成功能：调试前后比赛日获取函数' paragraph" }
if ($idx42 -eq 0) { throw "could not find the trailing '下一步工作' paragraph" }

# --- Edit 1: strip pPr/rPr from the "完成功能：调试前后比赛日获取函数" paragraph ---
$p40Range = $d.Paragraphs($idx40).Range

$frag40 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:ind w:left="1260" w:hangingChars="600" w:hanging="1260"/>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>完成功能</w:t>
</w:r>
<w:r>
<w:t>：</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>调试前后比赛日</w:t>
</w:r>
<w:r>
<w:t>获取函数</w:t>
</w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$p40Range.InsertXML($frag40)

# Re-resolve the second target paragraph (index may have shifted if the
# first edit changed the paragraph count, though it should not here).
$idx42 = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $ptext = $d.Paragraphs($i).Range.Text
    if ($ptext.StartsWith("下一步") -and $ptext.Contains("尝试") -and $ptext.Contains("数据库以供网页和远程应用访问")) {
        $idx42 = $i
    }
}
if ($idx42 -eq 0) { throw "could not re-find the trailing '下一步工作' paragraph" }

# --- Edit 2: replace the final paragraph with the rewritten paragraph plus
# the five new paragraphs describing 17/6/24's work. ---
$p42Range = $d.Paragraphs($idx42).Range

$frag42 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:ind w:left="1260" w:hangingChars="600" w:hanging="1260"/>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>下一步</w:t>
</w:r>
<w:r>
<w:t>工作：</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>尝试</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>将</w:t>
</w:r>
<w:r>
<w:t>数据</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>存入</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>sql</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>数据库以供网页和远程应用访问。</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:ind w:left="1260" w:hangingChars="600" w:hanging="1260"/>
</w:pPr>
</w:p>
<w:p>
<w:pPr>
<w:ind w:left="1260" w:hangingChars="600" w:hanging="1260"/>
</w:pPr>
<w:r>
<w:t>17/6/24</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:ind w:left="1260" w:hangingChars="600" w:hanging="1260"/>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>完成功能：</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>建立</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>MySQL</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>的</w:t>
</w:r>
<w:r>
<w:t>新用户，</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>建立</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>MySQL</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>里</w:t>
</w:r>
<w:r>
<w:t>的爬虫数据库并转换成</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>unicode</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>字符集，建立</w:t>
</w:r>
<w:r>
<w:t>表，</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>将队伍数据</w:t>
</w:r>
<w:r>
<w:t>存入</w:t>
</w:r>
<w:r>
<w:t>MN</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>ySQL</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>。</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:ind w:left="1260" w:hangingChars="600" w:hanging="1260"/>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>已解决问题：</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>将队伍数据存入</w:t>
</w:r>
<w:r>
<w:t>M</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>ySQL</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>。</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>MySQL</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>中表</w:t>
</w:r>
<w:r>
<w:t>的</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>某</w:t>
</w:r>
<w:r>
<w:t>一列</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>即使</w:t>
</w:r>
<w:r>
<w:t>是</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>INT</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>，</w:t>
</w:r>
<w:r>
<w:t>用</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>python</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>插入</w:t>
</w:r>
<w:r>
<w:t>的时候，</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>命令中</w:t>
</w:r>
<w:r>
<w:t>也应该用</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>\</w:t>
</w:r>
<w:r>
<w:t>”%s\”</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>作为</w:t>
</w:r>
<w:r>
<w:t>占位符。</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:ind w:left="1260" w:hangingChars="600" w:hanging="1260"/>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>未解决问题</w:t>
</w:r>
<w:r>
<w:t>：</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>由于</w:t>
</w:r>
<w:r>
<w:t>添加了</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>MySQL</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>的</w:t>
</w:r>
<w:r>
<w:t>部分，需要</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>为不同使用者或</w:t>
</w:r>
<w:r>
<w:t>不同抓取主机初始化</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>MySQL</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>的指导</w:t>
</w:r>
<w:r>
<w:t>或者输出</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>以</w:t>
</w:r>
<w:r>
<w:t>初始化的</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>MySQL</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>爬虫数据库</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>以供导入，或者</w:t>
</w:r>
<w:r>
<w:t>需要在</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>github</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>上</w:t>
</w:r>
<w:r>
<w:t>建立</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>readme</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>，</w:t>
</w:r>
<w:r>
<w:t>介绍</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>MySQL</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>中</w:t>
</w:r>
<w:r>
<w:t>的表结构。</w:t>
</w:r>
</w:p>
<w:p>
<w:pPr>
<w:ind w:left="1260" w:hangingChars="600" w:hanging="1260"/>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
</w:pPr>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>下一步工作：</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>尝试</w:t>
</w:r>
<w:r>
<w:t>将</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>本轮</w:t>
</w:r>
<w:r>
<w:t>、上轮、下轮比赛结果</w:t>
</w:r>
<w:r>
<w:rPr>
<w:rFonts w:hint="eastAsia"/>
</w:rPr>
<w:t>存入</w:t>
</w:r>
<w:r>
<w:t>数据库。</w:t>
</w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$p42Range.InsertXML($frag42)
